# Update the "想去人数" (want-to-go count) column F values on the
# "展览" and "全部类型" worksheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4285
$ws1.Range("F3").Value = 2430
$ws1.Range("F4").Value = 481
$ws1.Range("F8").Value = 213
$ws1.Range("F9").Value = 126
$ws1.Range("F10").Value = 134
$ws1.Range("F11").Value = 153
$ws1.Range("F12").Value = 1592
$ws1.Range("F13").Value = 293
$ws1.Range("F14").Value = 3320
$ws1.Range("F15").Value = 224

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4285
$ws4.Range("F3").Value = 2430
$ws4.Range("F4").Value = 481
$ws4.Range("F10").Value = 213
$ws4.Range("F11").Value = 126
$ws4.Range("F12").Value = 134
$ws4.Range("F13").Value = 153
$ws4.Range("F16").Value = 1592
$ws4.Range("F17").Value = 293
$ws4.Range("F18").Value = 3320
$ws4.Range("F19").Value = 224
